# Add help >>> Spec (BMW,OBC,DCB,5DH,NISSAN) information dialog
# The author clicked through every sheet tab (leaving a new cursor position
# behind on each one) and, while on the NISSAN_OBC sheet, replaced the old
# single-row "NISSAN" placeholder with three new contact-pair rows.

$wb = $excel.ActiveWorkbook

# --- BMW --------------------------------------------------------------
$ws = $wb.Worksheets.Item("BMW")
[void]$ws.Activate()
[void]$ws.Range("H32").Select()

# --- DAI_OBC ------------------------------------------------------------
$ws = $wb.Worksheets.Item("DAI_OBC")
[void]$ws.Activate()
[void]$ws.Range("I25").Select()

# --- DAI_DCB1.2 -----------------------------------------------------------
$ws = $wb.Worksheets.Item("DAI_DCB1.2")
[void]$ws.Activate()
[void]$ws.Range("I23").Select()

# --- DAI_DCB1.2H ----------------------------------------------------------
$ws = $wb.Worksheets.Item("DAI_DCB1.2H")
[void]$ws.Activate()
[void]$ws.Range("H25").Select()

# --- DAI_DCB2.0 -----------------------------------------------------------
$ws = $wb.Worksheets.Item("DAI_DCB2.0")
[void]$ws.Activate()
[void]$ws.Range("I35").Select()

# --- REN_5DH --------------------------------------------------------------
$ws = $wb.Worksheets.Item("REN_5DH")
[void]$ws.Activate()
[void]$ws.Range("H31").Select()

# --- NISSAN_OBC -------------------------------------------------------
# Replace the placeholder "NISSAN" row with real contact-pair data and
# add two more rows underneath.
$ws = $wb.Worksheets.Item("NISSAN_OBC")
[void]$ws.Activate()

$ws.Range("A2").Value = "MPE to MP1"
$ws.Range("A3").Value = "MPE to MP2"
$ws.Range("B3").Value = 10
$ws.Range("A4").Value = "MPR to MP3"
$ws.Range("B4").Value = 10

[void]$ws.Range("L13").Select()

# --- CUSTOM -----------------------------------------------------------
$ws = $wb.Worksheets.Item("CUSTOM")
[void]$ws.Activate()
[void]$ws.Range("L28").Select()

# Leave NISSAN_OBC as the active sheet/tab, matching activeTab="6" in the
# saved workbook.
$ws = $wb.Worksheets.Item("NISSAN_OBC")
[void]$ws.Activate()
[void]$ws.Range("L13").Select()
